# Update Excel files for new format
# - Add a new "Setup" sheet (inserted before the existing sheet) holding the
#   per-conflict-pair timing parameters (column A mirrors the row headers of
#   the conflict matrix, column B holds the gap/time values).
# - Rename the original sheet ("Sheet1") to "ConflictMatrix".
# - Restore sensible selections on both sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Setup" sheet before the existing (only) sheet ----------
# Worksheets.Add() with no arguments inserts the new sheet immediately before
# the currently active sheet - exactly like it was done in Excel.
$setup = $wb.Worksheets.Add()
$setup.Name = "Setup"

# --- 2. Rename the pre-existing sheet to "ConflictMatrix" -------------------
$matrix = $wb.Worksheets.Item("Sheet1")
$matrix.Name = "ConflictMatrix"

# --- 3. Populate the Setup sheet --------------------------------------------
# Column A: the conflict id / timing key (same values as column A of the
# conflict matrix). Column B: the associated gap time (seconds).
$values = @(
    @(1.1, 1.5),
    @(2.1, 1.5),
    @(5.1, 1.5),
    @(6.1, 1.5),
    @(7.1, 1.5),
    @(8.1, 1.5),
    @(9.1, 1.5),
    @(10.1, 1.5),
    @(11.1, 1.5),
    @(12.1, 1.5),
    @(86.1, 1.5),
    @(35.1, 1.5),
    @(26.1, 1.5),
    @(36.2, 1.5),
    @(88.1, 1.5),
    @(37.2, 1.5),
    @(28.1, 1.5),
    @(38.2, 1.5),
    @(31.2, 1.5),
    @(22, 1.5),
    @(32.2, 1.5),
    @(35.2, 1.5),
    @(36.1, 1.5),
    @(37.1, 1.5),
    @(38.1, 1.5),
    @(31.1, 1.5),
    @(32.1, 1.5),
    @(42, 30)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $setup.Cells.Item($row, 1).Value = $values[$i][0]
    $setup.Cells.Item($row, 2).Value = $values[$i][1]
}

# Give the Setup column A the same visual formatting (fill / font / border /
# number format) used for the matching rows in the conflict matrix header
# column, by copying the formats across.
$matrix.Range("A1:A28").Copy()
$setup.Range("A1:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Selections / active views -------------------------------------------
# ConflictMatrix: keep its zoom, select A1:A28, and it is no longer the
# selected tab.
$matrix.Range("A1:A28").Select()

# Setup: becomes the active/selected tab, with B29 selected (first empty row
# beneath the table).
$setup.Range("B29").Select()
$setup.Activate()
